$d = $word.ActiveDocument

# Replace the block of paragraphs from "You will be creating BMI calculator..."
# through the paragraph that originally held the drawing, with the new
# reordered/merged content (drawing moved inline into the first paragraph,
# "means" split into "me" + drawing + "ans", remaining paragraphs moved up,
# and the drawing's old paragraph left empty).
$startPara = $d.Paragraphs(6)
$endPara = $d.Paragraphs(10)
$r = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3141BCD2" w14:textId="77777777" w:rsidR="002546CC" w:rsidRDefault="002546CC" w:rsidP="00A70489"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr><w:t>You will be creating BMI calculator for this exam. BMI me</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="79A2CD1D" wp14:editId="07EF3A34"><wp:extent cx="5869576" cy="1607502"/><wp:effectExtent l="0" t="0" r="0" b="5715"/><wp:docPr id="802339805" name="Picture 1"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="802339805" name="Picture 1"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId8"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="5935708" cy="1625614"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr><w:t>ans</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve"> Body Mass Index value for human beings given weight in pounds and height in inches. </w:t></w:r></w:p><w:p w14:paraId="41F59510" w14:textId="03FF3359" w:rsidR="00A70489" w:rsidRDefault="00A70489" w:rsidP="00A70489"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r w:rsidRPr="002F5296"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">Create a project in </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr><w:t>X</w:t></w:r><w:r w:rsidRPr="002F5296"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr><w:t>Code with the name &#8220;</w:t></w:r><w:r w:rsidRPr="002F5296"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/><w:b/></w:rPr><w:t>LastName</w:t></w:r><w:r w:rsidRPr="002546CC"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/><w:b/></w:rPr><w:t>_</w:t></w:r><w:r w:rsidR="00F25833"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/><w:b/></w:rPr><w:t>Practice</w:t></w:r><w:r w:rsidR="0049164B"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/><w:b/></w:rPr><w:t>Exam01</w:t></w:r><w:r w:rsidRPr="002F5296"><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr><w:t xml:space="preserve">&#8221;. </w:t></w:r></w:p><w:p w14:paraId="05645950" w14:textId="77777777" w:rsidR="009625CC" w:rsidRPr="00F25833" w:rsidRDefault="009625CC" w:rsidP="00A70489"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/><w:bCs/></w:rPr><w:t>While creating the app make sure to follow minimum deployment and project format as shown below:</w:t></w:r></w:p><w:p w14:paraId="3207C14D" w14:textId="77777777" w:rsidR="00F25833" w:rsidRPr="00F25833" w:rsidRDefault="00F25833" w:rsidP="00F25833"><w:pPr><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr></w:pPr></w:p><w:p w14:paraId="75C30503" w14:textId="77888D98" w:rsidR="002546CC" w:rsidRDefault="0091571A" w:rsidP="009625CC"><w:pPr><w:ind w:left="720"/><w:rPr><w:rFonts w:ascii="Cambria" w:hAnsi="Cambria" w:cstheme="minorHAnsi"/></w:rPr></w:pPr></w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)

Write-Output "Block replaced."
